{"js": "// Update the worksheet date and the 25 three-digit-by-one-digit division\n// problems/answers. Every \"old\" string below occurs exactly once in the\n// document, so a scoped search + whole-match replace is unambiguous.\nconst replacements = [\n  [\"2025-09-06 Saturday\", \"2025-09-07 Sunday\"],\n  [\"305\u00f79=33, 8\", \"760\u00f77=108, 4\"],\n  [\"698\u00f72=349, 0\", \"683\u00f74=170, 3\"],\n  [\"498\u00f73=166, 0\", \"163\u00f78=20, 3\"],\n  [\"654\u00f74=163, 2\", \"664\u00f78=83, 0\"],\n  [\"878\u00f76=146, 2\", \"196\u00f79=21, 7\"],\n  [\"883\u00f74=220, 3\", \"980\u00f74=245, 0\"],\n  [\"134\u00f79=14, 8\", \"164\u00f78=20, 4\"],\n  [\"117\u00f73=39, 0\", \"873\u00f73=291, 0\"],\n  [\"934\u00f74=233, 2\", \"746\u00f72=373, 0\"],\n  [\"655\u00f75=131, 0\", \"390\u00f79=43, 3\"],\n  [\"167\u00f73=55, 2\", \"570\u00f74=142, 2\"],\n  [\"484\u00f78=60, 4\", \"959\u00f79=106, 5\"],\n  [\"490\u00f74=122, 2\", \"927\u00f77=132, 3\"],\n  [\"302\u00f74=75, 2\", \"358\u00f76=59, 4\"],\n  [\"918\u00f73=306, 0\", \"828\u00f74=207, 0\"],\n  [\"870\u00f79=96, 6\", \"741\u00f78=92, 5\"],\n  [\"500\u00f72=250, 0\", \"432\u00f78=54, 0\"],\n  [\"791\u00f77=113, 0\", \"534\u00f77=76, 2\"],\n  [\"957\u00f74=239, 1\", \"820\u00f74=205, 0\"],\n  [\"533\u00f74=133, 1\", \"542\u00f75=108, 2\"],\n  [\"332\u00f76=55, 2\", \"910\u00f73=303, 1\"],\n  [\"288\u00f72=144, 0\", \"486\u00f72=243, 0\"],\n  [\"599\u00f77=85, 4\", \"872\u00f76=145, 2\"],\n  [\"170\u00f78=21, 2\", \"710\u00f76=118, 2\"],\n  [\"159\u00f73=53, 0\", \"545\u00f73=181, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 three-digit-by-one-digit division\n# problems/answers using Word's Find & Replace (Execute) on the document's\n# full content range. Every \"old\" string occurs exactly once in the\n# document, so wdReplaceAll (with MatchCase) replaces exactly one instance\n# each, unambiguously.\n\n$d = $word.ActiveDocument\n\nfunction Replace-OneText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #          MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n\nReplace-OneText \"2025-09-06 Saturday\" \"2025-09-07 Sunday\"\n\nReplace-OneText \"305\u00f79=33, 8\" \"760\u00f77=108, 4\"\nReplace-OneText \"698\u00f72=349, 0\" \"683\u00f74=170, 3\"\nReplace-OneText \"498\u00f73=166, 0\" \"163\u00f78=20, 3\"\nReplace-OneText \"654\u00f74=163, 2\" \"664\u00f78=83, 0\"\nReplace-OneText \"878\u00f76=146, 2\" \"196\u00f79=21, 7\"\n\nReplace-OneText \"883\u00f74=220, 3\" \"980\u00f74=245, 0\"\nReplace-OneText \"134\u00f79=14, 8\" \"164\u00f78=20, 4\"\nReplace-OneText \"117\u00f73=39, 0\" \"873\u00f73=291, 0\"\nReplace-OneText \"934\u00f74=233, 2\" \"746\u00f72=373, 0\"\nReplace-OneText \"655\u00f75=131, 0\" \"390\u00f79=43, 3\"\n\nReplace-OneText \"167\u00f73=55, 2\" \"570\u00f74=142, 2\"\nReplace-OneText \"484\u00f78=60, 4\" \"959\u00f79=106, 5\"\nReplace-OneText \"490\u00f74=122, 2\" \"927\u00f77=132, 3\"\nReplace-OneText \"302\u00f74=75, 2\" \"358\u00f76=59, 4\"\nReplace-OneText \"918\u00f73=306, 0\" \"828\u00f74=207, 0\"\n\nReplace-OneText \"870\u00f79=96, 6\" \"741\u00f78=92, 5\"\nReplace-OneText \"500\u00f72=250, 0\" \"432\u00f78=54, 0\"\nReplace-OneText \"791\u00f77=113, 0\" \"534\u00f77=76, 2\"\nReplace-OneText \"957\u00f74=239, 1\" \"820\u00f74=205, 0\"\nReplace-OneText \"533\u00f74=133, 1\" \"542\u00f75=108, 2\"\n\nReplace-OneText \"332\u00f76=55, 2\" \"910\u00f73=303, 1\"\nReplace-OneText \"288\u00f72=144, 0\" \"486\u00f72=243, 0\"\nReplace-OneText \"599\u00f77=85, 4\" \"872\u00f76=145, 2\"\nReplace-OneText \"170\u00f78=21, 2\" \"710\u00f76=118, 2\"\nReplace-OneText \"159\u00f73=53, 0\" \"545\u00f73=181, 2\"\n"}
